$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "89.431.50"
$ws.Range("E2").Value = "  +11.25%  "

# Row 3
$ws.Range("D3").Value = "3.378.86"
$ws.Range("E3").Value = "  +6.06%  "

# Row 4
$ws.Range("E4").Value = "  -0.16%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "223.47"
$ws.Range("E5").Value = "  +6.30%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "645.96"
$ws.Range("E6").Value = "  +2.89%  "

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.339"
$ws.Range("E7").Value = "  +23.90%  "

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.999"
$ws.Range("E8").Value = "  -0.13%  "

# Row 9
$ws.Range("E9").Value = "  +4.97%  "

# Row 10
$ws.Range("D10").Value = "3.380.22"
$ws.Range("E10").Value = "  +6.17%  "

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.603"
$ws.Range("E11").Value = "  +2.06%  "

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.0000278"
$ws.Range("E12").Value = "  +7.57%  "

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.168"
$ws.Range("E13").Value = "  +2.29%  "

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "35.49"
$ws.Range("E14").Value = "  +11.03%  "

# Row 15
$ws.Range("D15").Value = "3.984.32"
$ws.Range("E15").Value = "  +5.55%  "

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "5.46"
$ws.Range("E16").Value = "  +3.48%  "

# Row 17
$ws.Range("D17").Value = "88.815.73"
$ws.Range("E17").Value = "  +10.37%  "

# Row 18
$ws.Range("D18").Value = "3.351.72"
$ws.Range("E18").Value = "  +4.93%  "

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "14.78"
$ws.Range("E19").Value = "  +3.71%  "

# Row 20
$ws.Range("E20").Value = "  +6.63%  "

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "474.94"
$ws.Range("E21").Value = "  +8.42%  "

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "9.22"
$ws.Range("E22").Value = "  +0.59%  "

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "5.47"
$ws.Range("E23").Value = "  +4.94%  "

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "13.88"
$ws.Range("E24").Value = "  +27.11%  "

# Row 25
$ws.Range("E25").Value = "  +7.06%  "

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "5.45"
$ws.Range("E26").Value = "  +15.98%  "

# Row 27
$ws.Range("E27").Value = "  +4.67%  "

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "79.60"
$ws.Range("E28").Value = "  +4.37%  "

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.216"
$ws.Range("E29").Value = "  +75.33%  "

# Row 30
$ws.Range("E30").Value = "  +6.43%  "

# Row 31
$ws.Range("E31").Value = "  -0.13%  "

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "603.07"
$ws.Range("E32").Value = "  +7.79%  "

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "9.38"
$ws.Range("E33").Value = "  +4.81%  "

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.994"
$ws.Range("E34").Value = "  -0.33%  "

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.56"
$ws.Range("E35").Value = "  +7.75%  "

# Row 36
$ws.Range("E36").Value = "  +3.62%  "

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.152"
$ws.Range("E37").Value = "  +0.78%  "

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "24.16"
$ws.Range("E38").Value = "  +4.75%  "

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "6.90"
$ws.Range("E39").Value = "  +22.29%  "

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.424"
$ws.Range("E40").Value = "  +4.06%  "

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.998"
$ws.Range("E41").Value = "  -0.12%  "

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "21.75"
$ws.Range("E42").Value = "  +4.74%  "

# Row 43
$ws.Range("E43").Value = "  +15.94%  "

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "3.08"
$ws.Range("E44").Value = "  +13.10%  "

# Row 45
$ws.Range("B45").Value = "Aave"
$ws.Range("C45").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "193.45"
$ws.Range("E45").Value = "  +2.16%  "

# Row 46
$ws.Range("B46").Value = "USDe"
$ws.Range("C46").Value = "https://coinranking.com/coin/exbfr2U-0+usde-usde"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "1.00"
$ws.Range("E46").Value = "  +0.02%  "

# Row 47
$ws.Range("B47").Value = "Monero"
$ws.Range("C47").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "157.15"
$ws.Range("E47").Value = "  -3.64%  "

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "47.65"
$ws.Range("E48").Value = "  +11.26%  "

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.38"
$ws.Range("E49").Value = "  +6.29%  "

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.800"
$ws.Range("E50").Value = "  +2.14%  "

# Row 51
$ws.Range("B51").Value = "ARBITRUM"
$ws.Range("C51").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.667"
$ws.Range("E51").Value = "  +6.33%  "
